# "Generate Report for Archive"
# - Update status text "Ready for handoff" -> "In Translation" on the
#   Overview sheet (columns E/F, row 2) and on each language sheet
#   (column C, row 2).
# - Narrow the now-shorter "In Translation" status columns:
#     Overview!E:F and zh-cn!C / de-de!C.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Column widths were narrowed to fit the shorter status text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
